$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2079.2
$ws.Range("I29").Value = 599
$ws.Range("K29").Value = 1797
$ws.Range("M29").Value = -1516
$ws.Range("H38").Value = 571.6667
$ws.Range("I38").Value = 571.6667
$ws.Range("K38").Value = 1715.0001
$ws.Range("M38").Value = -1343.0001
$ws.Range("H39").Value = 75.933334
$ws.Range("I39").Value = 45.642857
$ws.Range("K39").Value = 136.928571
$ws.Range("M39").Value = 159.071429
$ws.Range("H58").Value = 1898.8334
$ws.Range("I58").Value = 1864.3334
$ws.Range("K58").Value = 5593.0002
$ws.Range("M58").Value = -5443.0002
$ws.Range("H62").Value = 130178.375
$ws.Range("I62").Value = 170747
$ws.Range("K62").Value = 170747
$ws.Range("M62").Value = -170123
$ws.Range("H65").Value = 130178.375
$ws.Range("I65").Value = 170747
$ws.Range("K65").Value = 853735
$ws.Range("M65").Value = -850615
$ws.Range("H135").Value = 963.5294
$ws.Range("I135").Value = 797.26666
$ws.Range("K135").Value = 7175.39994
$ws.Range("M135").Value = -4640.39994
$ws.Range("H137").Value = 3370.0576
$ws.Range("J137").Value = 11535.7
$ws.Range("L137").Value = 34607.10000000001
$ws.Range("N137").Value = -39707.10000000001
$ws.Range("H141").Value = 50789.473
$ws.Range("I141").Value = 53083.332
$ws.Range("J141").Value = 9500
$ws.Range("K141").Value = 159249.996
$ws.Range("L141").Value = 28500
$ws.Range("M141").Value = -154069.996
$ws.Range("N141").Value = -38860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 2943.75
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 2943.75
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 2943.75
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -3691.75
$ws.Range("H32").Value = 25467.512
$ws.Range("I32").Value = 14696.395
$ws.Range("K32").Value = 14696.395
$ws.Range("M32").Value = -14409.395
$ws.Range("H61").Value = 1141
$ws.Range("I61").Value = 767.0345
$ws.Range("J61").Value = 2690.2856
$ws.Range("K61").Value = 767.0345
$ws.Range("L61").Value = 2690.2856
$ws.Range("M61").Value = -555.0345
$ws.Range("N61").Value = -3114.2856
$ws.Range("H74").Value = 1455.96
$ws.Range("I74").Value = 1321.1305
$ws.Range("J74").Value = 3006.5
$ws.Range("K74").Value = 1321.1305
$ws.Range("L74").Value = 3006.5
$ws.Range("M74").Value = -447.1305
$ws.Range("N74").Value = -4754.5
$ws.Range("H77").Value = 1455.96
$ws.Range("I77").Value = 1321.1305
$ws.Range("J77").Value = 3006.5
$ws.Range("K77").Value = 6605.6525
$ws.Range("L77").Value = 15032.5
$ws.Range("M77").Value = -2237.6525
$ws.Range("N77").Value = -23768.5
$ws.Range("H122").Value = 2333.6924
$ws.Range("I122").Value = 2196.3333
$ws.Range("J122").Value = 2642.75
$ws.Range("K122").Value = 6588.999899999999
$ws.Range("L122").Value = 7928.25
$ws.Range("M122").Value = -4138.999899999999
$ws.Range("N122").Value = -12828.25
$ws.Range("H128").Value = 77999
$ws.Range("J128").Value = 77999
$ws.Range("L128").Value = 77999
$ws.Range("N128").Value = -87959
$ws.Range("H136").Value = 1141
$ws.Range("I136").Value = 767.0345
$ws.Range("J136").Value = 2690.2856
$ws.Range("K136").Value = 2301.1035
$ws.Range("L136").Value = 8070.8568
$ws.Range("M136").Value = 248.8964999999998
$ws.Range("N136").Value = -13170.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 921
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 921
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 921
$ws.Range("M29").Value = ""
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H134").Value = 4017.8
$ws.Range("I134").Value = 3947.6428
$ws.Range("K134").Value = 11842.9284
$ws.Range("M134").Value = -9307.928400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 306.75
$ws.Range("I5").Value = 306.75
$ws.Range("K5").Value = 920.25
$ws.Range("M5").Value = -808.25
$ws.Range("H75").Value = 7382.6113
$ws.Range("I75").Value = 4225.5
$ws.Range("J75").Value = 8284.643
$ws.Range("K75").Value = 12676.5
$ws.Range("L75").Value = 24853.929
$ws.Range("M75").Value = -11678.5
$ws.Range("N75").Value = -26849.929
$ws.Range("H78").Value = 7382.6113
$ws.Range("I78").Value = 4225.5
$ws.Range("J78").Value = 8284.643
$ws.Range("K78").Value = 38029.5
$ws.Range("L78").Value = 74561.787
$ws.Range("M78").Value = -33037.5
$ws.Range("N78").Value = -84545.787
$ws.Range("H87").Value = 10330.429
$ws.Range("I87").Value = 10330.429
$ws.Range("K87").Value = 30991.287
$ws.Range("M87").Value = -29743.287
$ws.Range("H90").Value = 10330.429
$ws.Range("I90").Value = 10330.429
$ws.Range("K90").Value = 92973.861
$ws.Range("M90").Value = -86733.861
$ws.Range("H107").Value = 963.8
$ws.Range("I107").Value = 1029
$ws.Range("K107").Value = 3087
$ws.Range("M107").Value = -1167
$ws.Range("H117").Value = 3608.3076
$ws.Range("J117").Value = 3831.2173
$ws.Range("L117").Value = 11493.6519
$ws.Range("N117").Value = -18377.6519
$ws.Range("H132").Value = 1750
$ws.Range("J132").Value = 1900
$ws.Range("L132").Value = 17100
$ws.Range("N132").Value = -22160
$ws.Range("H135").Value = 306.75
$ws.Range("I135").Value = 306.75
$ws.Range("K135").Value = 2760.75
$ws.Range("M135").Value = -225.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 9999
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 9999
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -10575
$ws.Range("H40").Value = 18672
$ws.Range("I40").Value = 13008
$ws.Range("K40").Value = 13008
$ws.Range("M40").Value = -12857
$ws.Range("H58").Value = 24247.25
$ws.Range("I58").Value = 22000
$ws.Range("K58").Value = 22000
$ws.Range("M58").Value = -21723
$ws.Range("H80").Value = 5621
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 6328
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 6328
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -8324
$ws.Range("H83").Value = 5621
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 6328
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 31640
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -41624
$ws.Range("H97").Value = 22271.516
$ws.Range("I97").Value = 33492.45
$ws.Range("J97").Value = 1869.8182
$ws.Range("K97").Value = 33492.45
$ws.Range("L97").Value = 1869.8182
$ws.Range("M97").Value = -32996.45
$ws.Range("N97").Value = -2861.8182
$ws.Range("H132").Value = 1983.931
$ws.Range("I132").Value = 1847.7391
$ws.Range("J132").Value = 2506
$ws.Range("K132").Value = 5543.2173
$ws.Range("L132").Value = 7518
$ws.Range("M132").Value = -3013.2173
$ws.Range("N132").Value = -12578

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 15000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = ""
$ws.Range("H133").Value = 99992.664
$ws.Range("J133").Value = 99992.664
$ws.Range("L133").Value = 99992.664
$ws.Range("N133").Value = -105052.664
$ws.Range("H136").Value = 2638.963
$ws.Range("I136").Value = 2116.2273
$ws.Range("K136").Value = 6348.6819
$ws.Range("M136").Value = -3798.6819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1501600
$ws.Range("J14").Value = 2133.3333
$ws.Range("L14").Value = 2133.3333
$ws.Range("N14").Value = -2469.3333
$ws.Range("H81").Value = 5171.909
$ws.Range("I81").Value = 5171.909
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 10343.818
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -9282.817999999999
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 5171.909
$ws.Range("I84").Value = 5171.909
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 51719.09
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -46415.09
$ws.Range("N84").Value = ""
$ws.Range("H131").Value = 97914
$ws.Range("J131").Value = 97914
$ws.Range("L131").Value = 97914
$ws.Range("N131").Value = -107994
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = ""
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
